$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.009538531303406
$ws.Range("B1").Value = 2.120188236236572
$ws.Range("C1").Value = 6.307126522064209
$ws.Range("D1").Value = 1.507797956466675
$ws.Range("E1").Value = 1.349408745765686
